$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 6: "Baseline 2010 C109" ----
$ws.Range("A6").Value = "CW3M"
$ws.Range("B6").Value = "Baseline 2010 C109"
$ws.Range("C6").Value = 2010

$ws.Range("D6").Value = 1143.6110839999999
$ws.Range("D6").NumberFormat = "0.00"
$ws.Range("E6").Value = 1921.3682859999999
$ws.Range("E6").NumberFormat = "0.00"
$ws.Range("F6").Value = 4.0370929999999996
$ws.Range("F6").NumberFormat = "0.00"
$ws.Range("G6").Value = 197.01855499999999
$ws.Range("G6").NumberFormat = "0.00"
$ws.Range("H6").Value = 73.459366000000003
$ws.Range("H6").NumberFormat = "0.00"

$ws.Range("I6").Value = 105.737251
$ws.Range("I6").NumberFormat = "0.00"
$ws.Range("I6").Interior.Color = 65535

$ws.Range("J6").Value = 62.789425000000001
$ws.Range("J6").NumberFormat = "0.00"
$ws.Range("K6").Value = 717.21026600000005
$ws.Range("K6").NumberFormat = "0.00"
$ws.Range("L6").Value = 75.975761000000006
$ws.Range("L6").NumberFormat = "0.00"

$ws.Range("M6").Value = 982.77160600000002
$ws.Range("M6").NumberFormat = "0.00"
$ws.Range("M6").Interior.Color = 65535

$ws.Range("N6").Value = 1219.040649
$ws.Range("N6").NumberFormat = "0.00"

$ws.Range("O6").Value = 6231.375
$ws.Range("O6").NumberFormat = "0"
$ws.Range("P6").Value = 162867.046875
$ws.Range("P6").NumberFormat = "0"

$ws.Range("Q6").Value = -387.443939
$ws.Range("Q6").NumberFormat = "0.00"
$ws.Range("Q6").Interior.Color = 65535

$ws.Range("R6").Value = -0.112458
$ws.Range("R6").NumberFormat = "0.000000"
$ws.Range("R6").Interior.Color = 65535

$ws.Range("S6").Value = 2010

# ---- Row 7: "C109 with no springs" ----
$ws.Range("B7").Value = "C109 with no springs"
$ws.Range("C7").Value = 2010

$ws.Range("D7").Value = 1143.6110839999999
$ws.Range("D7").NumberFormat = "0.00"
$ws.Range("E7").Value = 1921.3682859999999
$ws.Range("E7").NumberFormat = "0.00"
$ws.Range("F7").Value = 4.0370929999999996
$ws.Range("F7").NumberFormat = "0.00"
$ws.Range("G7").Value = 0
$ws.Range("G7").NumberFormat = "0.00"
$ws.Range("H7").Value = 73.459366000000003
$ws.Range("H7").NumberFormat = "0.00"
$ws.Range("I7").Value = 4.2416460000000002
$ws.Range("I7").NumberFormat = "0.00"
$ws.Range("J7").Value = 62.789425000000001
$ws.Range("J7").NumberFormat = "0.00"
$ws.Range("K7").Value = 717.21569799999997
$ws.Range("K7").NumberFormat = "0.00"
$ws.Range("L7").Value = 75.975761000000006
$ws.Range("L7").NumberFormat = "0.00"
$ws.Range("M7").Value = 1078.0303960000001
$ws.Range("M7").NumberFormat = "0.00"
$ws.Range("N7").Value = 1219.301514
$ws.Range("N7").NumberFormat = "0.00"

$ws.Range("O7").Value = 6231.375
$ws.Range("O7").NumberFormat = "0"
$ws.Range("P7").Value = 162867.046875
$ws.Range("P7").NumberFormat = "0"

$ws.Range("Q7").Value = 6.5953189999999999
$ws.Range("Q7").NumberFormat = "0.00"

$ws.Range("R7").Value = 0.0020960000000000002
$ws.Range("R7").NumberFormat = "0.000000"

$ws.Range("S7").Value = 2010

# ---- Sheet/window state to match target ----
$null = $ws.Range("B8").Select()
